$d = $word.ActiveDocument

# Locate the paragraph that starts the section ("What is a corporation?" /
# style SectionHeadnote) and the paragraph that ends the block we are
# collapsing (the final SectionHeadnote, "This is the second chapter...").
$startPara = $null
$endPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Style.NameLocal -eq "Section Headnote") {
        if ($startPara -eq $null) {
            $startPara = $p
        }
        $endPara = $p
    }
}

# Build the merged headnote text: section/resource numbers, titles and
# case-text bodies are concatenated together; the (now unused) resource /
# section headnote bodies are dropped.
$target = "1.1" + `
    "Case of the District Number 1" + `
    "This is the body of case 1." + `
    "1.2" + `
    "Case of the District Number 2" + `
    "highlighted: content to highlight; elided: content to elide; replaced: content to replace; commented: content to comment; highlighted2: second highlight content;`n" + `
    "2" + `
    "Section Two"

# Replace the first paragraph's text in place (keeps it a single run /
# single paragraph so xml:space="preserve" round-trips correctly), then
# delete everything from there through the end of the last paragraph that
# is being collapsed into it.
$startPara.Range.Text = $target

$delRange = $d.Range($startPara.Range.End, $endPara.Range.End)
$delRange.Delete()
